# Scheduled market-data refresh for the Leve-profit workbook.
# The runner re-pulls currentAveragePrice* (H:L) from the market board and
# recomputes LeveProfitNQ/HQ (M:N). When a leve's NQ or HQ reward quantity is
# zero, that side has no meaningful profit figure, so its M/N cell is cleared
# instead of left holding a stale number (and the other side gains one).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 266.14285
$ws.Range("I2").Value = 310
$ws.Range("K2").Value = 310
$ws.Range("M2").Value = -197
$ws.Range("H40").Value = 1211.5883
$ws.Range("I40").Value = 1170.6207
$ws.Range("K40").Value = 1170.6207
$ws.Range("M40").Value = -995.6206999999999
$ws.Range("H43").Value = 11255.95
$ws.Range("I43").Value = 15055.857
$ws.Range("J43").Value = 9209.846
$ws.Range("K43").Value = 15055.857
$ws.Range("L43").Value = 9209.846
$ws.Range("M43").Value = -14986.857
$ws.Range("N43").Value = -9347.846
$ws.Range("H58").Value = 184.2
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H64").Value = 4399.6665
$ws.Range("I64").Value = 3999
$ws.Range("K64").Value = 3999
$ws.Range("M64").Value = -3751
$ws.Range("H67").Value = 4399.6665
$ws.Range("I67").Value = 3999
$ws.Range("K67").Value = 3999
$ws.Range("M67").Value = -3141
$ws.Range("H127").Value = 145473.42
$ws.Range("J127").Value = 337372.34
$ws.Range("L127").Value = 1012117.02
$ws.Range("N127").Value = -1022037.02

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 21502
$ws.Range("I19").Value = 21502
$ws.Range("K19").Value = 21502
$ws.Range("M19").Value = -21273
$ws.Range("H37").Value = 30000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 30000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 30000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -30546
$ws.Range("H61").Value = 4337.933
$ws.Range("I61").Value = 3928.3845
$ws.Range("K61").Value = 3928.3845
$ws.Range("M61").Value = -3716.3845
$ws.Range("H74").Value = 64049.188
$ws.Range("I74").Value = 84691.25
$ws.Range("K74").Value = 84691.25
$ws.Range("M74").Value = -83817.25
$ws.Range("H77").Value = 64049.188
$ws.Range("I77").Value = 84691.25
$ws.Range("K77").Value = 423456.25
$ws.Range("M77").Value = -419088.25
$ws.Range("H132").Value = 3439.1365
$ws.Range("I132").Value = 3735.8667
$ws.Range("K132").Value = 11207.6001
$ws.Range("M132").Value = -8677.6001
$ws.Range("H136").Value = 4337.933
$ws.Range("I136").Value = 3928.3845
$ws.Range("K136").Value = 11785.1535
$ws.Range("M136").Value = -9235.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11047.952
$ws.Range("I86").Value = 8401.207
$ws.Range("J86").Value = 16952.23
$ws.Range("K86").Value = 8401.207
$ws.Range("L86").Value = 16952.23
$ws.Range("M86").Value = -7278.207
$ws.Range("N86").Value = -19198.23
$ws.Range("H89").Value = 11047.952
$ws.Range("I89").Value = 8401.207
$ws.Range("J89").Value = 16952.23
$ws.Range("K89").Value = 42006.035
$ws.Range("L89").Value = 84761.14999999999
$ws.Range("M89").Value = -36390.035
$ws.Range("N89").Value = -95993.14999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 323.85715
$ws.Range("J22").Value = 316.75
$ws.Range("L22").Value = 316.75
$ws.Range("N22").Value = -1016.75
$ws.Range("H41").Value = 10999.3125
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 10999.3125
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 10999.3125
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -11855.3125
$ws.Range("H47").Value = 19300
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 19300
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 19300
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -20432
$ws.Range("H62").Value = 5247.3335
$ws.Range("I62").Value = 5748.5
$ws.Range("K62").Value = 5748.5
$ws.Range("M62").Value = -5124.5
$ws.Range("H65").Value = 5247.3335
$ws.Range("I65").Value = 5748.5
$ws.Range("K65").Value = 28742.5
$ws.Range("M65").Value = -25622.5
$ws.Range("H105").Value = 1673.75
$ws.Range("I105").Value = 1521.5385
$ws.Range("K105").Value = 1521.5385
$ws.Range("M105").Value = 225.4614999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 25641550
$ws.Range("I107").Value = 33333616
$ws.Range("J107").Value = 1326.3334
$ws.Range("K107").Value = 100000848
$ws.Range("L107").Value = 3979.0002
$ws.Range("M107").Value = -99998928
$ws.Range("N107").Value = -7819.0002
$ws.Range("H113").Value = 499.33334
$ws.Range("I113").Value = 499.33334
$ws.Range("K113").Value = 1498.00002
$ws.Range("M113").Value = 671.9999800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 270.5
$ws.Range("I2").Value = 208.61539
$ws.Range("J2").Value = 385.42856
$ws.Range("K2").Value = 208.61539
$ws.Range("L2").Value = 385.42856
$ws.Range("M2").Value = -95.61538999999999
$ws.Range("N2").Value = -611.4285600000001
$ws.Range("H57").Value = 10591
$ws.Range("H80").Value = 2604.4167
$ws.Range("I80").Value = 2569.8
$ws.Range("J80").Value = 2629.1428
$ws.Range("K80").Value = 2569.8
$ws.Range("L80").Value = 2629.1428
$ws.Range("M80").Value = -1571.8
$ws.Range("N80").Value = -4625.1428
$ws.Range("H83").Value = 2604.4167
$ws.Range("I83").Value = 2569.8
$ws.Range("J83").Value = 2629.1428
$ws.Range("K83").Value = 12849
$ws.Range("L83").Value = 13145.714
$ws.Range("M83").Value = -7857
$ws.Range("N83").Value = -23129.714
$ws.Range("H122").Value = 1378.3
$ws.Range("I122").Value = 1122.875
$ws.Range("K122").Value = 3368.625
$ws.Range("M122").Value = -918.625
$ws.Range("H126").Value = 24068
$ws.Range("I126").Value = 37935.555
$ws.Range("J126").Value = 3266.6667
$ws.Range("K126").Value = 113806.665
$ws.Range("L126").Value = 9800.000100000001
$ws.Range("M126").Value = -111336.665
$ws.Range("N126").Value = -14740.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 2000
$ws.Range("I11").Value = 2000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1860
$ws.Range("N11").ClearContents()
$ws.Range("H22").Value = 5946.077
$ws.Range("J22").Value = 10616.5
$ws.Range("L22").Value = 10616.5
$ws.Range("N22").Value = -11206.5
$ws.Range("H27").Value = 5946.077
$ws.Range("J27").Value = 10616.5
$ws.Range("L27").Value = 10616.5
$ws.Range("N27").Value = -10830.5
$ws.Range("H55").Value = 308.33334
$ws.Range("I55").Value = 279.16666
$ws.Range("K55").Value = 279.16666
$ws.Range("M55").Value = -106.16666
$ws.Range("H82").Value = 2157.5715
$ws.Range("H85").Value = 2157.5715
$ws.Range("H127").Value = 56943.332
$ws.Range("J127").Value = 56943.332
$ws.Range("L127").Value = 56943.332
$ws.Range("N127").Value = -66863.33199999999
$ws.Range("H132").Value = 3652.9707
$ws.Range("I132").Value = 2933.5
$ws.Range("J132").Value = 5379.7
$ws.Range("K132").Value = 8800.5
$ws.Range("L132").Value = 16139.1
$ws.Range("M132").Value = -6270.5
$ws.Range("N132").Value = -21199.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2478.7727
$ws.Range("I132").Value = 2406.5715
$ws.Range("K132").Value = 7219.7145
$ws.Range("M132").Value = -4689.7145
$ws.Range("H136").Value = 1541.0968
$ws.Range("I136").Value = 1251
$ws.Range("J136").Value = 3499.25
$ws.Range("K136").Value = 3753
$ws.Range("L136").Value = 10497.75
$ws.Range("M136").Value = -1203
$ws.Range("N136").Value = -15597.75
